$wb = $excel.ActiveWorkbook

# Sheet ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1948769.2
$ws.Range("I33").Value = 300
$ws.Range("K33").Value = 300
$ws.Range("M33").Value = -71

# Sheet ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 835499
$ws.Range("I39").Value = 1058244.9
$ws.Range("J39").Value = 202
$ws.Range("K39").Value = 3174734.7
$ws.Range("L39").Value = 606
$ws.Range("M39").Value = -3174438.7
$ws.Range("N39").Value = -1198

# Sheet ARM row 9
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 12000
$ws.Range("J9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("N9").Value = -12340

# Sheet ARM row 20
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 12000
$ws.Range("J20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("N20").Value = -12540

# Sheet ARM row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 80006
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# Sheet ARM row 29
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10616

# Sheet ARM row 33
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 20000
$ws.Range("I33").Value = 20000
$ws.Range("K33").Value = 20000
$ws.Range("M33").Value = -19671

# Sheet ARM row 38
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 15000
$ws.Range("J38").Value = 15000
$ws.Range("L38").Value = 15000
$ws.Range("N38").Value = -15934

# Sheet ARM row 39
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 16
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Sheet ARM row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1790
$ws.Range("I41").Value = 1790
$ws.Range("K41").Value = 1790
$ws.Range("M41").Value = -1376

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3148.3333
$ws.Range("I122").Value = 2722.1428
$ws.Range("J122").Value = 4640
$ws.Range("K122").Value = 8166.428400000001
$ws.Range("L122").Value = 13920
$ws.Range("M122").Value = -5716.428400000001
$ws.Range("N122").Value = -18820

# Sheet ARM row 34
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 49800
$ws.Range("J34").Value = 49800
$ws.Range("L34").Value = 49800
$ws.Range("N34").Value = -50028

# Sheet BSM row 39
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 1000
$ws.Range("J39").Value = 1000
$ws.Range("L39").Value = 1000
$ws.Range("N39").Value = -1778

# Sheet BSM row 46
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Sheet BSM row 49
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 21849.75
$ws.Range("J49").Value = 21849.75
$ws.Range("L49").Value = 21849.75
$ws.Range("N49").Value = -22327.75

# Sheet BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 60238.094
$ws.Range("J132").Value = 60238.094
$ws.Range("L132").Value = 60238.094
$ws.Range("N132").Value = -70358.094

# Sheet CRP row 29
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Sheet CRP row 35
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 2499.875
$ws.Range("I35").Value = 1142.7142
$ws.Range("J35").Value = 12000
$ws.Range("K35").Value = 1142.7142
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = -848.7141999999999
$ws.Range("N35").Value = -12588

# Sheet CRP row 42
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 9166.666999999999
$ws.Range("I42").Value = 3750
$ws.Range("J42").Value = 20000
$ws.Range("K42").Value = 3750
$ws.Range("L42").Value = 20000
$ws.Range("M42").Value = -3157
$ws.Range("N42").Value = -21186

# Sheet CRP row 55
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 13229.2
$ws.Range("I55").Value = 10382
$ws.Range("J55").Value = 17500
$ws.Range("K55").Value = 10382
$ws.Range("L55").Value = 17500
$ws.Range("M55").Value = -10067
$ws.Range("N55").Value = -18130

# Sheet CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 910.7143
$ws.Range("I26").Value = 220.14285
$ws.Range("J26").Value = 1256
$ws.Range("K26").Value = 660.4285500000001
$ws.Range("L26").Value = 3768
$ws.Range("M26").Value = -372.4285500000001
$ws.Range("N26").Value = -4344

# Sheet CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 100476.6
$ws.Range("I46").Value = 294.33334
$ws.Range("J46").Value = 250750
$ws.Range("K46").Value = 883.0000200000001
$ws.Range("L46").Value = 752250
$ws.Range("M46").Value = -792.0000200000001
$ws.Range("N46").Value = -752432

# Sheet CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 10226.846
$ws.Range("J55").Value = 3139.9
$ws.Range("L55").Value = 9419.700000000001
$ws.Range("N55").Value = -9773.700000000001

# Sheet GSM row 29
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 13674.333
$ws.Range("I29").Value = 1007
$ws.Range("J29").Value = 20008
$ws.Range("K29").Value = 1007
$ws.Range("L29").Value = 20008
$ws.Range("M29").Value = -717
$ws.Range("N29").Value = -20588

# Sheet GSM row 4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 10000
$ws.Range("K4").Value = 10000
$ws.Range("M4").Value = -9887

# Sheet LTW row 5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 22390
$ws.Range("J5").Value = 22390
$ws.Range("L5").Value = 22390
$ws.Range("N5").Value = -22616

# Sheet LTW row 25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 47002.668
$ws.Range("I25").Value = 70004
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 70004
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = -69774
$ws.Range("N25").Value = -1460

# Sheet LTW row 28
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 10000
$ws.Range("I28").Value = 10000
$ws.Range("K28").Value = 10000
$ws.Range("M28").Value = -9768

# Sheet LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 6046
$ws.Range("I32").Value = 565
$ws.Range("J32").Value = 9700
$ws.Range("K32").Value = 565
$ws.Range("L32").Value = 9700
$ws.Range("M32").Value = -248
$ws.Range("N32").Value = -10334

# Sheet LTW row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 2666.5
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 3333
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 3333
$ws.Range("M34").Value = -1828
$ws.Range("N34").Value = -3677

# Sheet LTW row 37
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 10000
$ws.Range("I37").Value = 10000
$ws.Range("K37").Value = 10000
$ws.Range("M37").Value = -9893

# Sheet LTW row 41
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 6915.222
$ws.Range("J41").Value = 6915.222
$ws.Range("L41").Value = 6915.222
$ws.Range("N41").Value = -7791.222

# Sheet LTW row 47
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 12343.571
$ws.Range("J47").Value = 12343.571
$ws.Range("L47").Value = 12343.571
$ws.Range("N47").Value = -13323.571

# Sheet LTW row 52
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 12343.571
$ws.Range("J52").Value = 12343.571
$ws.Range("L52").Value = 12343.571
$ws.Range("N52").Value = -12809.571

# Sheet LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1414.8286
$ws.Range("I55").Value = 1860.1111
$ws.Range("K55").Value = 1860.1111
$ws.Range("M55").Value = -1687.1111

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4403.846
$ws.Range("I132").Value = 6919.636
$ws.Range("J132").Value = 2558.9333
$ws.Range("K132").Value = 20758.908
$ws.Range("L132").Value = 7676.7999
$ws.Range("M132").Value = -18228.908
$ws.Range("N132").Value = -12736.7999

# Sheet WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Sheet WVR row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 9125
$ws.Range("J20").Value = 9125
$ws.Range("L20").Value = 9125
$ws.Range("N20").Value = -9605

# Sheet WVR row 22
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Sheet WVR row 28
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 9400
$ws.Range("J28").Value = 9450
$ws.Range("L28").Value = 9450
$ws.Range("N28").Value = -10146

# Sheet WVR row 29
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 23021.8
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 28652.25
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 28652.25
$ws.Range("M29").Value = -210
$ws.Range("N29").Value = -29232.25

# Sheet WVR row 30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 32504.5
$ws.Range("I30").Value = 32504.5
$ws.Range("K30").Value = 32504.5
$ws.Range("M30").Value = -32397.5

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2349.558
$ws.Range("I132").Value = 2561.9
$ws.Range("J132").Value = 1859.5385
$ws.Range("K132").Value = 7685.700000000001
$ws.Range("L132").Value = 5578.6155
$ws.Range("M132").Value = -5155.700000000001
$ws.Range("N132").Value = -10638.6155
